$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2580
$ws.Range("F5").Value = 914
$ws.Range("F7").Value = 1869
$ws.Range("F8").Value = 1725
$ws.Range("F9").Value = 192
$ws.Range("F11").Value = 2397
$ws.Range("F12").Value = 507
$ws.Range("F13").Value = 181
$ws.Range("F17").Value = 100
$ws.Range("F18").Value = 8752
$ws.Range("F20").Value = 6829
$ws.Range("F21").Value = 11032
$ws.Range("F23").Value = 189
$ws.Range("F24").Value = 222
$ws.Range("F25").Value = 304
$ws.Range("F26").Value = 530
$ws.Range("F27").Value = 2436
$ws.Range("F30").Value = 2274
$ws.Range("F31").Value = 388
$ws.Range("F32").Value = 28
$ws.Range("F33").Value = 4459
$ws.Range("F34").Value = 597
$ws.Range("F35").Value = 233
$ws.Range("F36").Value = 19
$ws.Range("F37").Value = 456

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 10

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 790
$ws.Range("F5").Value = 73

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 790
$ws.Range("F5").Value = 73
$ws.Range("F7").Value = 2580
$ws.Range("F9").Value = 914
$ws.Range("F11").Value = 1869
$ws.Range("F13").Value = 1725
$ws.Range("F15").Value = 192
$ws.Range("F17").Value = 507
$ws.Range("F18").Value = 181
$ws.Range("F22").Value = 100
$ws.Range("F23").Value = 8752
$ws.Range("F25").Value = 6829
$ws.Range("F26").Value = 11032
$ws.Range("F27").Value = 10
$ws.Range("F29").Value = 189
$ws.Range("F30").Value = 222
$ws.Range("F31").Value = 304
$ws.Range("F33").Value = 530
$ws.Range("F37").Value = 28
$ws.Range("F38").Value = 4459
$ws.Range("F45").Value = 456
